$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the part numbers
$ws.Range("B2").Value = "DS2310LGWHT-LF"
$ws.Range("B3").Value = "DS3410LGWHT-LF"

# Update the quantities
$ws.Range("C2").Value = 3.0
$ws.Range("C3").Value = 1.0

# Remove the now-obsolete row 4
$ws.Range("B4:C4").Delete()
